# Commit: "commit add them ki nang song"
#
# The existing first data row (STT=1, TH Hoà Bình, Lớp 1, 12, 14/07/2018, ...)
# is removed. This shifts the remaining two data rows up by one (old row 8
# becomes row 7, old row 9 becomes row 8) and the STT numbering is fixed up
# (1, 2). The event date-range subtitle is also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old first data row (row 7); Excel shifts everything below it up.
$ws.Rows(7).Delete()

# Renumber the STT column for the two remaining data rows.
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2

# Update the subtitle date range.
$ws.Range("A4").Value = "Từ ngày 18/07/2018 tới ngày 18/08/2018"

# Update the "SỐ LƯỢNG" (quantity) for the second remaining row.
$ws.Range("D8").Value = 80
